$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value2 = 7000.6
$ws.Range("J116").Value2 = 7750.75
$ws.Range("L116").Value2 = 7750.75
$ws.Range("N116").Value2 = -14634.75

$ws.Range("H127").Value2 = 1351.1666
$ws.Range("J127").Value2 = 1785.6666
$ws.Range("L127").Value2 = 5356.9998
$ws.Range("N127").Value2 = -15276.9998

$ws.Range("H132").Value2 = 4424
$ws.Range("I132").Value2 = 4363.0527
$ws.Range("K132").Value2 = 13089.1581
$ws.Range("M132").Value2 = -10559.1581

$ws.Range("H137").Value2 = 2343
$ws.Range("I137").Value2 = 1941.5834
$ws.Range("J137").Value2 = 4751.5
$ws.Range("K137").Value2 = 5824.7502
$ws.Range("L137").Value2 = 14254.5
$ws.Range("M137").Value2 = -3274.7502
$ws.Range("N137").Value2 = -19354.5

$ws.Range("H141").Value2 = 1980.16
$ws.Range("I141").Value2 = 1422.7
$ws.Range("J141").Value2 = 4210
$ws.Range("K141").Value2 = 4268.1
$ws.Range("L141").Value2 = 12630
$ws.Range("M141").Value2 = 911.8999999999996
$ws.Range("N141").Value2 = -22990

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 3472.4167
$ws.Range("I61").Value2 = 3401.7273
$ws.Range("K61").Value2 = 3401.7273
$ws.Range("M61").Value2 = -3189.7273

$ws.Range("H122").Value2 = 3774.9167
$ws.Range("I122").Value2 = 2912.5
$ws.Range("J122").Value2 = 5499.75
$ws.Range("K122").Value2 = 8737.5
$ws.Range("L122").Value2 = 16499.25
$ws.Range("M122").Value2 = -6287.5
$ws.Range("N122").Value2 = -21399.25

$ws.Range("H132").Value2 = 23632.625
$ws.Range("I132").Value2 = 2642
$ws.Range("J132").Value2 = 65613.875
$ws.Range("K132").Value2 = 7926
$ws.Range("L132").Value2 = 196841.625
$ws.Range("M132").Value2 = -5396
$ws.Range("N132").Value2 = -201901.625

$ws.Range("H136").Value2 = 3472.4167
$ws.Range("I136").Value2 = 3401.7273
$ws.Range("K136").Value2 = 10205.1819
$ws.Range("M136").Value2 = -7655.1819

$ws.Range("H138").Value2 = 48183.2
$ws.Range("J138").Value2 = 48183.2
$ws.Range("L138").Value2 = 48183.2
$ws.Range("N138").Value2 = -58463.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value2 = 15674.223
$ws.Range("J81").Value2 = 15674.223
$ws.Range("L81").Value2 = 15674.223
$ws.Range("N81").Value2 = -17796.223

$ws.Range("H84").Value2 = 15674.223
$ws.Range("J84").Value2 = 15674.223
$ws.Range("L84").Value2 = 47022.669
$ws.Range("N84").Value2 = -57630.669

$ws.Range("H99").Value2 = 1421.6154
$ws.Range("I99").Value2 = 1377
$ws.Range("K99").Value2 = 1377
$ws.Range("M99").Value2 = 121

$ws.Range("H107").Value2 = 970
$ws.Range("I107").Value2 = 837.5
$ws.Range("K107").Value2 = 837.5
$ws.Range("M107").Value2 = 1082.5

$ws.Range("H134").Value2 = 3452.7188
$ws.Range("I134").Value2 = 3531.8386
$ws.Range("K134").Value2 = 10595.5158
$ws.Range("M134").Value2 = -8060.515800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2600.8604
$ws.Range("I31").Value2 = 1078
$ws.Range("J31").Value2 = 3336.0344
$ws.Range("K31").Value2 = 1078
$ws.Range("L31").Value2 = 3336.0344
$ws.Range("M31").Value2 = -783
$ws.Range("N31").Value2 = -3926.0344

$ws.Range("H34").Value2 = 2600.8604
$ws.Range("I34").Value2 = 1078
$ws.Range("J34").Value2 = 3336.0344
$ws.Range("K34").Value2 = 1078
$ws.Range("L34").Value2 = 3336.0344
$ws.Range("M34").Value2 = -876
$ws.Range("N34").Value2 = -3740.0344

$ws.Range("H53").Value2 = 34893
$ws.Range("J53").Value2 = 34893
$ws.Range("L53").Value2 = 34893
$ws.Range("N53").Value2 = -36107

$ws.Range("H107").Value2 = 1836.2941
$ws.Range("I107").Value2 = 1578.1666
$ws.Range("K107").Value2 = 1578.1666
$ws.Range("M107").Value2 = 341.8334

$ws.Range("H134").Value2 = 1251.3914
$ws.Range("I134").Value2 = 1186.7778
$ws.Range("J134").Value2 = 1484
$ws.Range("K134").Value2 = 3560.3334
$ws.Range("L134").Value2 = 4452
$ws.Range("M134").Value2 = -1025.3334
$ws.Range("N134").Value2 = -9522

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value2 = 4509.6875
$ws.Range("I109").Value2 = 909
$ws.Range("J109").Value2 = 5340.615
$ws.Range("K109").Value2 = 2727
$ws.Range("L109").Value2 = 16021.845
$ws.Range("M109").Value2 = -1687
$ws.Range("N109").Value2 = -18101.845

$ws.Range("H131").Value2 = 711
$ws.Range("J131").Value2 = 722.2105
$ws.Range("L131").Value2 = 2166.6315
$ws.Range("N131").Value2 = -12246.6315

$ws.Range("H140").Value2 = 2536.5217
$ws.Range("I140").Value2 = 1289.2307
$ws.Range("J140").Value2 = 4158
$ws.Range("K140").Value2 = 3867.6921
$ws.Range("L140").Value2 = 12474
$ws.Range("M140").Value2 = 1312.3079
$ws.Range("N140").Value2 = -22834

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value2 = 13000.167
$ws.Range("J92").Value2 = 13000.167
$ws.Range("L92").Value2 = 13000.167
$ws.Range("N92").Value2 = -16744.167

$ws.Range("H102").Value2 = 2724.6667
$ws.Range("I102").Value2 = 2883.7646
$ws.Range("J102").Value2 = 2338.2856
$ws.Range("K102").Value2 = 2883.7646
$ws.Range("L102").Value2 = 2338.2856
$ws.Range("M102").Value2 = -1261.7646
$ws.Range("N102").Value2 = -5582.2856

$ws.Range("H132").Value2 = 21193.654
$ws.Range("I132").Value2 = 1412.5714
$ws.Range("K132").Value2 = 4237.7142
$ws.Range("M132").Value2 = -1707.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 3720.15
$ws.Range("I40").Value2 = 3282.5454
$ws.Range("J40").Value2 = 4255
$ws.Range("K40").Value2 = 3282.5454
$ws.Range("L40").Value2 = 4255
$ws.Range("M40").Value2 = -3146.5454
$ws.Range("N40").Value2 = -4527

$ws.Range("H55").Value2 = 244.77777
$ws.Range("I55").Value2 = 190.2
$ws.Range("J55").Value2 = 313
$ws.Range("K55").Value2 = 190.2
$ws.Range("L55").Value2 = 313
$ws.Range("M55").Value2 = -17.19999999999999
$ws.Range("N55").Value2 = -659

$ws.Range("H122").Value2 = 1228416.1
$ws.Range("I122").Value2 = 1785150.8
$ws.Range("K122").Value2 = 5355452.4
$ws.Range("M122").Value2 = -5353002.4

$ws.Range("H132").Value2 = 2489.3333
$ws.Range("I132").Value2 = 2054
$ws.Range("J132").Value2 = 4666
$ws.Range("K132").Value2 = 6162
$ws.Range("L132").Value2 = 13998
$ws.Range("M132").Value2 = -3632
$ws.Range("N132").Value2 = -19058

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value2 = 1246.1538
$ws.Range("I81").Value2 = 1287.625
$ws.Range("J81").Value2 = 1179.8
$ws.Range("K81").Value2 = 2575.25
$ws.Range("L81").Value2 = 2359.6
$ws.Range("M81").Value2 = -1514.25
$ws.Range("N81").Value2 = -4481.6

$ws.Range("H84").Value2 = 1246.1538
$ws.Range("I84").Value2 = 1287.625
$ws.Range("J84").Value2 = 1179.8
$ws.Range("K84").Value2 = 12876.25
$ws.Range("L84").Value2 = 11798
$ws.Range("M84").Value2 = -7572.25
$ws.Range("N84").Value2 = -22406

$ws.Range("H92").Value2 = 32275
$ws.Range("J92").Value2 = 32275
$ws.Range("L92").Value2 = 32275
$ws.Range("N92").Value2 = -37267

$ws.Range("H113").Value2 = 3001.8333
$ws.Range("I113").Value2 = 3546.2
$ws.Range("K113").Value2 = 10638.6
$ws.Range("M113").Value2 = -8468.599999999999

$ws.Range("H122").Value2 = 1245.2
$ws.Range("I122").Value2 = 1317
$ws.Range("J122").Value2 = 1137.5
$ws.Range("K122").Value2 = 3951
$ws.Range("L122").Value2 = 3412.5
$ws.Range("M122").Value2 = -1501
$ws.Range("N122").Value2 = -8312.5

$ws.Range("H136").Value2 = 31253338
$ws.Range("I136").Value2 = 41668156
$ws.Range("K136").Value2 = 125004468
$ws.Range("M136").Value2 = -125001918

$ws.Range("H140").Value2 = 38782.25
$ws.Range("J140").Value2 = 38782.25
$ws.Range("L140").Value2 = 38782.25
$ws.Range("N140").Value2 = -49142.25

$ws.Range("H141").Value2 = 80000
$ws.Range("J141").Value2 = 80000
$ws.Range("L141").Value2 = 80000
$ws.Range("N141").Value2 = -90360
